$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Daily update: append today's row (33) with the new values, then carry
# the "last row" date-format style down from row 32 to row 33, and
# restore row 32 to the regular (non-last-row) date format.

$ws.Range("A33").Value = 45617
$ws.Range("B33").Value = 82
$ws.Range("C33").Value = 68
$ws.Range("D33").Value = 79

$ws.Range("A33").NumberFormat = $ws.Range("A32").NumberFormat
$ws.Range("A32").NumberFormat = $ws.Range("A31").NumberFormat
